# Generate Report for Handback
$wb = $excel.ActiveWorkbook

# Row 3 on every sheet (Overview, zh-cn, de-de) corresponds to the
# 77c7dea1-... file whose Status was "Ready for handoff" and is now
# "Handback transform failed" (shared string reused across sheets).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# New Error Detail (column K) for row 3 on the zh-cn and de-de sheets.
$wsZhCn.Range("K3").Value = "Handback file name: xdo4kk1g.qtc is different with handoff file name: 77c7dea1-77ae-4509-8510-f7ea3e20695b.fa32823acaf01bd51dd51b8072be1b05ae320786.zh-cn."
$wsDeDe.Range("K3").Value = "Handback file name: xdo4kk1g.qtc is different with handoff file name: 77c7dea1-77ae-4509-8510-f7ea3e20695b.fa32823acaf01bd51dd51b8072be1b05ae320786.de-de."
